# Auto-generated edit script: apply numeric updates from the commit diff.
# Each worksheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) gets its changed
# cells in the "currentAveragePrice.../LeveProfit..." columns (H-N) updated
# to match the refreshed market-price snapshot.

$wb = $excel.ActiveWorkbook


# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 67.78570999999999
$ws.Range("I11").Value = 67.78570999999999
$ws.Range("K11").Value = 67.78570999999999
$ws.Range("M11").Value = 72.21429000000001
$ws.Range("H28").Value = 2361.2
$ws.Range("I28").Value = 2236.3333
$ws.Range("K28").Value = 2236.3333
$ws.Range("M28").Value = -1751.3333
$ws.Range("H33").Value = 412.27777
$ws.Range("I33").Value = 417.8125
$ws.Range("J33").Value = 368
$ws.Range("K33").Value = 417.8125
$ws.Range("L33").Value = 368
$ws.Range("M33").Value = -188.8125
$ws.Range("N33").Value = -826
$ws.Range("H43").Value = 3501.3333
$ws.Range("I43").Value = 3501
$ws.Range("J43").Value = 3502
$ws.Range("K43").Value = 3501
$ws.Range("L43").Value = 3502
$ws.Range("M43").Value = -3432
$ws.Range("N43").Value = -3640
$ws.Range("H88").Value = 2516.2942
$ws.Range("I88").Value = 2288.889
$ws.Range("J88").Value = 2772.125
$ws.Range("K88").Value = 2288.889
$ws.Range("L88").Value = 2772.125
$ws.Range("M88").Value = -1882.889
$ws.Range("N88").Value = -3584.125
$ws.Range("H91").Value = 2516.2942
$ws.Range("I91").Value = 2288.889
$ws.Range("J91").Value = 2772.125
$ws.Range("K91").Value = 2288.889
$ws.Range("L91").Value = 2772.125
$ws.Range("M91").Value = -884.8890000000001
$ws.Range("N91").Value = -5580.125
$ws.Range("H116").Value = 6847.909
$ws.Range("I116").Value = 5332.75
$ws.Range("K116").Value = 5332.75
$ws.Range("M116").Value = -1890.75
$ws.Range("H132").Value = 3676.4167
$ws.Range("I132").Value = 4211.7
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 12635.1
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -10105.1
$ws.Range("N132").Value = -8060
$ws.Range("H137").Value = 1199
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1300.8334
$ws.Range("I2").Value = 1501
$ws.Range("K2").Value = 1501
$ws.Range("M2").Value = -1388
$ws.Range("H45").Value = 5503.25
$ws.Range("I45").Value = 5333
$ws.Range("K45").Value = 5333
$ws.Range("M45").Value = -4956
$ws.Range("H61").Value = 3054.8
$ws.Range("I61").Value = 2582.8
$ws.Range("K61").Value = 2582.8
$ws.Range("M61").Value = -2370.8
$ws.Range("H74").Value = 3000
$ws.Range("I74").Value = 3000
$ws.Range("K74").Value = 3000
$ws.Range("M74").Value = -2126
$ws.Range("H77").Value = 3000
$ws.Range("I77").Value = 3000
$ws.Range("K77").Value = 15000
$ws.Range("M77").Value = -10632
$ws.Range("H80").Value = 44999
$ws.Range("J80").Value = 69999
$ws.Range("L80").Value = 69999
$ws.Range("N80").Value = -71995
$ws.Range("H83").Value = 44999
$ws.Range("J83").Value = 69999
$ws.Range("L83").Value = 209997
$ws.Range("N83").Value = -219981
$ws.Range("H110").Value = 4752.8335
$ws.Range("I110").Value = 5304.643
$ws.Range("K110").Value = 5304.643
$ws.Range("M110").Value = -3259.643
$ws.Range("H116").Value = 1300.8334
$ws.Range("I116").Value = 1501
$ws.Range("K116").Value = 1501
$ws.Range("M116").Value = 793
$ws.Range("H132").Value = 7582.5
$ws.Range("I132").Value = 7582.5
$ws.Range("K132").Value = 22747.5
$ws.Range("M132").Value = -20217.5
$ws.Range("H136").Value = 3054.8
$ws.Range("I136").Value = 2582.8
$ws.Range("K136").Value = 7748.400000000001
$ws.Range("M136").Value = -5198.400000000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1300.8334
$ws.Range("I3").Value = 1501
$ws.Range("K3").Value = 1501
$ws.Range("M3").Value = -1387
$ws.Range("H94").Value = 1280.375
$ws.Range("I94").Value = 2495.6667
$ws.Range("J94").Value = 999.9231
$ws.Range("K94").Value = 2495.6667
$ws.Range("L94").Value = 999.9231
$ws.Range("M94").Value = -2044.6667
$ws.Range("N94").Value = -1901.9231
$ws.Range("H105").Value = 4702.5
$ws.Range("I105").Value = 4702.5
$ws.Range("K105").Value = 4702.5
$ws.Range("M105").Value = -2955.5
$ws.Range("H138").Value = 146000
$ws.Range("J138").Value = 146000
$ws.Range("L138").Value = 146000
$ws.Range("N138").Value = -156280

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5927.778
$ws.Range("J31").Value = 4592.857
$ws.Range("L31").Value = 4592.857
$ws.Range("N31").Value = -5182.857
$ws.Range("H34").Value = 5927.778
$ws.Range("J34").Value = 4592.857
$ws.Range("L34").Value = 4592.857
$ws.Range("N34").Value = -4996.857
$ws.Range("H36").Value = 3000
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H40").Value = 3000
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H60").Value = 29799.6
$ws.Range("I60").Value = 29799.6
$ws.Range("K60").Value = 29799.6
$ws.Range("M60").Value = -29288.6
$ws.Range("H99").Value = 4702.4
$ws.Range("I99").Value = 5170.6665
$ws.Range("K99").Value = 5170.6665
$ws.Range("M99").Value = -3672.6665
$ws.Range("H105").Value = 1999
$ws.Range("I105").Value = 1999
$ws.Range("K105").Value = 1999
$ws.Range("M105").Value = -252
$ws.Range("H126").Value = 4702.4
$ws.Range("I126").Value = 5170.6665
$ws.Range("K126").Value = 15511.9995
$ws.Range("M126").Value = -13041.9995
$ws.Range("H134").Value = 6587.3335
$ws.Range("I134").Value = 6587.3335
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 19762.0005
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -17227.0005
$ws.Range("N134").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 6333.1665
$ws.Range("I49").Value = 7299.8
$ws.Range("K49").Value = 21899.4
$ws.Range("M49").Value = -21743.4
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H107").Value = 97.5
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
$ws.Range("H113").Value = 2120.1428
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2120.1428
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 6360.428400000001
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -10700.4284

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4622
$ws.Range("I122").Value = 3749.25
$ws.Range("J122").Value = 5494.75
$ws.Range("K122").Value = 11247.75
$ws.Range("L122").Value = 16484.25
$ws.Range("M122").Value = -8797.75
$ws.Range("N122").Value = -21384.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 2824.75
$ws.Range("I9").Value = 599.6667
$ws.Range("K9").Value = 599.6667
$ws.Range("M9").Value = -375.6667
$ws.Range("H35").Value = 1258
$ws.Range("I35").Value = 1094.8572
$ws.Range("K35").Value = 1094.8572
$ws.Range("M35").Value = -758.8571999999999
$ws.Range("H101").Value = 15787.333
$ws.Range("J101").Value = 15787.333
$ws.Range("L101").Value = 15787.333
$ws.Range("N101").Value = -22277.333
$ws.Range("H136").Value = 12752.723
$ws.Range("I136").Value = 12397.857
$ws.Range("J136").Value = 13994.75
$ws.Range("K136").Value = 37193.571
$ws.Range("L136").Value = 41984.25
$ws.Range("M136").Value = -34643.571
$ws.Range("N136").Value = -47084.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 3339699
$ws.Range("J4").Value = 3757036.2
$ws.Range("L4").Value = 3757036.2
$ws.Range("N4").Value = -3757262.2
$ws.Range("H33").Value = 292.5
$ws.Range("I33").Value = 292.5
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 292.5
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -42.5
$ws.Range("N33").ClearContents()
$ws.Range("H36").Value = 292.5
$ws.Range("I36").Value = 292.5
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 292.5
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -42.5
$ws.Range("N36").ClearContents()
$ws.Range("H96").Value = 4648.909
$ws.Range("I96").Value = 4363.2856
$ws.Range("K96").Value = 4363.2856
$ws.Range("M96").Value = -2990.2856
$ws.Range("H122").Value = 2866.2
$ws.Range("I122").Value = 2994
$ws.Range("J122").Value = 2568
$ws.Range("K122").Value = 8982
$ws.Range("L122").Value = 7704
$ws.Range("M122").Value = -6532
$ws.Range("N122").Value = -12604
